$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: new columns E (precio_num) and F (fecha_dia) ---
$ws.Cells.Item(1, 5).Value = "precio_num"
$ws.Cells.Item(1, 6).Value = "fecha_dia"

# --- Row 163: timestamp correction (A163) ---
$ws.Cells.Item(163, 1).Value = 45964.36344328704

# --- New rows 164-190 (A-D) ---
$ws.Cells.Item(164, 1).Value = 45966.36917886574
$ws.Cells.Item(164, 2).Value = "EVOWHEY PROTEIN"
$ws.Cells.Item(164, 3).Value = "2Kg"
$ws.Cells.Item(164, 4).Value = "37,90€"
$ws.Cells.Item(165, 1).Value = 45966.39759947917
$ws.Cells.Item(165, 2).Value = "EVOWHEY PROTEIN"
$ws.Cells.Item(165, 3).Value = "2Kg"
$ws.Cells.Item(165, 4).Value = "37,90€"
$ws.Cells.Item(166, 1).Value = 45966.45842030093
$ws.Cells.Item(166, 2).Value = "EVOWHEY PROTEIN"
$ws.Cells.Item(166, 3).Value = "2Kg"
$ws.Cells.Item(166, 4).Value = "37,90€"
$ws.Cells.Item(167, 1).Value = 45967.45853447916
$ws.Cells.Item(167, 2).Value = "EVOWHEY PROTEIN"
$ws.Cells.Item(167, 3).Value = "2Kg"
$ws.Cells.Item(167, 4).Value = "37,90€"
$ws.Cells.Item(168, 1).Value = 45968.45850428241
$ws.Cells.Item(168, 2).Value = "EVOWHEY PROTEIN"
$ws.Cells.Item(168, 3).Value = "2Kg"
$ws.Cells.Item(168, 4).Value = "37,90€"
$ws.Cells.Item(169, 1).Value = 45969.4584796412
$ws.Cells.Item(169, 2).Value = "EVOWHEY PROTEIN"
$ws.Cells.Item(169, 3).Value = "2Kg"
$ws.Cells.Item(169, 4).Value = "37,90€"
$ws.Cells.Item(170, 1).Value = 45970.45847402778
$ws.Cells.Item(170, 2).Value = "EVOWHEY PROTEIN"
$ws.Cells.Item(170, 3).Value = "2Kg"
$ws.Cells.Item(170, 4).Value = "37,90€"
$ws.Cells.Item(171, 1).Value = 45971.45847685185
$ws.Cells.Item(171, 2).Value = "EVOWHEY PROTEIN"
$ws.Cells.Item(171, 3).Value = "2Kg"
$ws.Cells.Item(171, 4).Value = "37,90€"
$ws.Cells.Item(172, 1).Value = 45972.4584791088
$ws.Cells.Item(172, 2).Value = "EVOWHEY PROTEIN"
$ws.Cells.Item(172, 3).Value = "2Kg"
$ws.Cells.Item(172, 4).Value = "31,96€"
$ws.Cells.Item(173, 1).Value = 45973.45851195602
$ws.Cells.Item(173, 2).Value = "EVOWHEY PROTEIN"
$ws.Cells.Item(173, 3).Value = "2Kg"
$ws.Cells.Item(173, 4).Value = "37,90€"
$ws.Cells.Item(174, 1).Value = 45974.45848630787
$ws.Cells.Item(174, 2).Value = "EVOWHEY PROTEIN"
$ws.Cells.Item(174, 3).Value = "2Kg"
$ws.Cells.Item(174, 4).Value = "37,90€"
$ws.Cells.Item(175, 1).Value = 45978.45848123843
$ws.Cells.Item(175, 2).Value = "EVOWHEY PROTEIN"
$ws.Cells.Item(175, 3).Value = "2Kg"
$ws.Cells.Item(175, 4).Value = "37,90€"
$ws.Cells.Item(176, 1).Value = 45979.45847527778
$ws.Cells.Item(176, 2).Value = "EVOWHEY PROTEIN"
$ws.Cells.Item(176, 3).Value = "2Kg"
$ws.Cells.Item(176, 4).Value = "37,90€"
$ws.Cells.Item(177, 1).Value = 45980.45847171296
$ws.Cells.Item(177, 2).Value = "EVOWHEY PROTEIN"
$ws.Cells.Item(177, 3).Value = "2Kg"
$ws.Cells.Item(177, 4).Value = "37,90€"
$ws.Cells.Item(178, 1).Value = 45981.45846827546
$ws.Cells.Item(178, 2).Value = "EVOWHEY PROTEIN"
$ws.Cells.Item(178, 3).Value = "2Kg"
$ws.Cells.Item(178, 4).Value = "37,90€"
$ws.Cells.Item(179, 1).Value = 45985.45850328704
$ws.Cells.Item(179, 2).Value = "EVOWHEY PROTEIN"
$ws.Cells.Item(179, 3).Value = "2Kg"
$ws.Cells.Item(179, 4).Value = "32,91€"
$ws.Cells.Item(180, 1).Value = 45986.39174679398
$ws.Cells.Item(180, 2).Value = "EVOWHEY PROTEIN"
$ws.Cells.Item(180, 3).Value = "2Kg"
$ws.Cells.Item(180, 4).Value = "32,91€"
$ws.Cells.Item(181, 1).Value = 45986.40505469908
$ws.Cells.Item(181, 2).Value = "EVOWHEY PROTEIN"
$ws.Cells.Item(181, 3).Value = "2Kg"
$ws.Cells.Item(181, 4).Value = "32,91€"
$ws.Cells.Item(182, 1).Value = 45986.40568177083
$ws.Cells.Item(182, 2).Value = "EVOWHEY PROTEIN"
$ws.Cells.Item(182, 3).Value = "2Kg"
$ws.Cells.Item(182, 4).Value = "32,91€"
$ws.Cells.Item(183, 1).Value = 45986.40632627315
$ws.Cells.Item(183, 2).Value = "EVOWHEY PROTEIN"
$ws.Cells.Item(183, 3).Value = "2Kg"
$ws.Cells.Item(183, 4).Value = "32,91€"
$ws.Cells.Item(184, 1).Value = 45986.40841292824
$ws.Cells.Item(184, 2).Value = "EVOWHEY PROTEIN"
$ws.Cells.Item(184, 3).Value = "2Kg"
$ws.Cells.Item(184, 4).Value = "32,91€"
$ws.Cells.Item(185, 1).Value = 45986.41840105324
$ws.Cells.Item(185, 2).Value = "EVOWHEY PROTEIN"
$ws.Cells.Item(185, 3).Value = "2Kg"
$ws.Cells.Item(185, 4).Value = "32,91€"
$ws.Cells.Item(186, 1).Value = 45986.41964681713
$ws.Cells.Item(186, 2).Value = "EVOWHEY PROTEIN"
$ws.Cells.Item(186, 3).Value = "2Kg"
$ws.Cells.Item(186, 4).Value = "32,91€"
$ws.Cells.Item(187, 1).Value = 45986.42271101852
$ws.Cells.Item(187, 2).Value = "EVOWHEY PROTEIN"
$ws.Cells.Item(187, 3).Value = "2Kg"
$ws.Cells.Item(187, 4).Value = "32,91€"
$ws.Cells.Item(188, 1).Value = 45986.42962280093
$ws.Cells.Item(188, 2).Value = "EVOWHEY PROTEIN"
$ws.Cells.Item(188, 3).Value = "2Kg"
$ws.Cells.Item(188, 4).Value = "32,91€"
$ws.Cells.Item(189, 1).Value = 45986.43288417824
$ws.Cells.Item(189, 2).Value = "EVOWHEY PROTEIN"
$ws.Cells.Item(189, 3).Value = "2Kg"
$ws.Cells.Item(189, 4).Value = "32,91€"
$ws.Cells.Item(190, 1).Value = 45986.43551783999
$ws.Cells.Item(190, 2).Value = "EVOWHEY PROTEIN"
$ws.Cells.Item(190, 3).Value = "2Kg"
$ws.Cells.Item(190, 4).Value = "32,91€"

# --- Column E (precio_num) and F (fecha_dia) for rows 2-189 (row 190 left blank per source) ---
$ws.Cells.Item(2, 5).Value = 37.9
$ws.Cells.Item(2, 6).Value = 45804
$ws.Cells.Item(3, 5).Value = 37.9
$ws.Cells.Item(3, 6).Value = 45804
$ws.Cells.Item(4, 5).Value = 37.9
$ws.Cells.Item(4, 6).Value = 45805
$ws.Cells.Item(5, 5).Value = 37.9
$ws.Cells.Item(5, 6).Value = 45806
$ws.Cells.Item(6, 5).Value = 37.9
$ws.Cells.Item(6, 6).Value = 45806
$ws.Cells.Item(7, 5).Value = 37.9
$ws.Cells.Item(7, 6).Value = 45807
$ws.Cells.Item(8, 5).Value = 37.9
$ws.Cells.Item(8, 6).Value = 45808
$ws.Cells.Item(9, 5).Value = 33.9
$ws.Cells.Item(9, 6).Value = 45809
$ws.Cells.Item(10, 5).Value = 34.9
$ws.Cells.Item(10, 6).Value = 45810
$ws.Cells.Item(11, 5).Value = 37.9
$ws.Cells.Item(11, 6).Value = 45811
$ws.Cells.Item(12, 5).Value = 34.9
$ws.Cells.Item(12, 6).Value = 45812
$ws.Cells.Item(13, 5).Value = 34.9
$ws.Cells.Item(13, 6).Value = 45813
$ws.Cells.Item(14, 5).Value = 37.9
$ws.Cells.Item(14, 6).Value = 45814
$ws.Cells.Item(15, 5).Value = 37.9
$ws.Cells.Item(15, 6).Value = 45815
$ws.Cells.Item(16, 5).Value = 34.9
$ws.Cells.Item(16, 6).Value = 45816
$ws.Cells.Item(17, 5).Value = 37.9
$ws.Cells.Item(17, 6).Value = 45817
$ws.Cells.Item(18, 5).Value = 37.9
$ws.Cells.Item(18, 6).Value = 45818
$ws.Cells.Item(19, 5).Value = 37.9
$ws.Cells.Item(19, 6).Value = 45833
$ws.Cells.Item(20, 5).Value = 37.9
$ws.Cells.Item(20, 6).Value = 45833
$ws.Cells.Item(21, 5).Value = 37.9
$ws.Cells.Item(21, 6).Value = 45833
$ws.Cells.Item(22, 5).Value = 37.9
$ws.Cells.Item(22, 6).Value = 45833
$ws.Cells.Item(23, 5).Value = 37.9
$ws.Cells.Item(23, 6).Value = 45833
$ws.Cells.Item(24, 5).Value = 37.9
$ws.Cells.Item(24, 6).Value = 45833
$ws.Cells.Item(25, 5).Value = 37.9
$ws.Cells.Item(25, 6).Value = 45833
$ws.Cells.Item(26, 5).Value = 37.9
$ws.Cells.Item(26, 6).Value = 45833
$ws.Cells.Item(27, 5).Value = 37.9
$ws.Cells.Item(27, 6).Value = 45833
$ws.Cells.Item(28, 5).Value = 37.9
$ws.Cells.Item(28, 6).Value = 45833
$ws.Cells.Item(29, 5).Value = 37.9
$ws.Cells.Item(29, 6).Value = 45833
$ws.Cells.Item(30, 5).Value = 37.9
$ws.Cells.Item(30, 6).Value = 45833
$ws.Cells.Item(31, 5).Value = 37.9
$ws.Cells.Item(31, 6).Value = 45833
$ws.Cells.Item(32, 5).Value = 37.9
$ws.Cells.Item(32, 6).Value = 45833
$ws.Cells.Item(33, 5).Value = 37.9
$ws.Cells.Item(33, 6).Value = 45833
$ws.Cells.Item(34, 5).Value = 37.9
$ws.Cells.Item(34, 6).Value = 45833
$ws.Cells.Item(35, 5).Value = 37.9
$ws.Cells.Item(35, 6).Value = 45833
$ws.Cells.Item(36, 5).Value = 37.9
$ws.Cells.Item(36, 6).Value = 45833
$ws.Cells.Item(37, 5).Value = 37.9
$ws.Cells.Item(37, 6).Value = 45833
$ws.Cells.Item(38, 5).Value = 37.9
$ws.Cells.Item(38, 6).Value = 45834
$ws.Cells.Item(39, 5).Value = 37.9
$ws.Cells.Item(39, 6).Value = 45853
$ws.Cells.Item(40, 5).Value = 37.9
$ws.Cells.Item(40, 6).Value = 45853
$ws.Cells.Item(41, 5).Value = 37.9
$ws.Cells.Item(41, 6).Value = 45853
$ws.Cells.Item(42, 5).Value = 37.9
$ws.Cells.Item(42, 6).Value = 45853
$ws.Cells.Item(43, 5).Value = 37.9
$ws.Cells.Item(43, 6).Value = 45853
$ws.Cells.Item(44, 5).Value = 37.9
$ws.Cells.Item(44, 6).Value = 45853
$ws.Cells.Item(45, 5).Value = 37.9
$ws.Cells.Item(45, 6).Value = 45853
$ws.Cells.Item(46, 5).Value = 37.9
$ws.Cells.Item(46, 6).Value = 45853
$ws.Cells.Item(47, 5).Value = 37.9
$ws.Cells.Item(47, 6).Value = 45853
$ws.Cells.Item(48, 5).Value = 37.9
$ws.Cells.Item(48, 6).Value = 45853
$ws.Cells.Item(49, 5).Value = 37.9
$ws.Cells.Item(49, 6).Value = 45853
$ws.Cells.Item(50, 5).Value = 37.9
$ws.Cells.Item(50, 6).Value = 45853
$ws.Cells.Item(51, 5).Value = 37.9
$ws.Cells.Item(51, 6).Value = 45853
$ws.Cells.Item(52, 5).Value = 37.9
$ws.Cells.Item(52, 6).Value = 45853
$ws.Cells.Item(53, 5).Value = 37.9
$ws.Cells.Item(53, 6).Value = 45853
$ws.Cells.Item(54, 5).Value = 37.9
$ws.Cells.Item(54, 6).Value = 45853
$ws.Cells.Item(55, 5).Value = 37.9
$ws.Cells.Item(55, 6).Value = 45853
$ws.Cells.Item(56, 5).Value = 37.9
$ws.Cells.Item(56, 6).Value = 45853
$ws.Cells.Item(57, 5).Value = 37.9
$ws.Cells.Item(57, 6).Value = 45853
$ws.Cells.Item(58, 5).Value = 37.9
$ws.Cells.Item(58, 6).Value = 45853
$ws.Cells.Item(59, 5).Value = 37.9
$ws.Cells.Item(59, 6).Value = 45854
$ws.Cells.Item(60, 5).Value = 37.9
$ws.Cells.Item(60, 6).Value = 45854
$ws.Cells.Item(61, 5).Value = 37.9
$ws.Cells.Item(61, 6).Value = 45855
$ws.Cells.Item(62, 5).Value = 37.9
$ws.Cells.Item(62, 6).Value = 45855
$ws.Cells.Item(63, 5).Value = 37.9
$ws.Cells.Item(63, 6).Value = 45855
$ws.Cells.Item(64, 5).Value = 37.9
$ws.Cells.Item(64, 6).Value = 45856
$ws.Cells.Item(65, 5).Value = 37.9
$ws.Cells.Item(65, 6).Value = 45857
$ws.Cells.Item(66, 5).Value = 37.9
$ws.Cells.Item(66, 6).Value = 45858
$ws.Cells.Item(67, 5).Value = 37.9
$ws.Cells.Item(67, 6).Value = 45859
$ws.Cells.Item(68, 5).Value = 37.9
$ws.Cells.Item(68, 6).Value = 45860
$ws.Cells.Item(69, 5).Value = 37.9
$ws.Cells.Item(69, 6).Value = 45861
$ws.Cells.Item(70, 5).Value = 37.9
$ws.Cells.Item(70, 6).Value = 45862
$ws.Cells.Item(71, 5).Value = 37.9
$ws.Cells.Item(71, 6).Value = 45863
$ws.Cells.Item(72, 5).Value = 37.9
$ws.Cells.Item(72, 6).Value = 45864
$ws.Cells.Item(73, 5).Value = 33.9
$ws.Cells.Item(73, 6).Value = 45865
$ws.Cells.Item(74, 5).Value = 33.9
$ws.Cells.Item(74, 6).Value = 45866
$ws.Cells.Item(75, 5).Value = 37.9
$ws.Cells.Item(75, 6).Value = 45867
$ws.Cells.Item(76, 5).Value = 37.9
$ws.Cells.Item(76, 6).Value = 45868
$ws.Cells.Item(77, 5).Value = 33.9
$ws.Cells.Item(77, 6).Value = 45869
$ws.Cells.Item(78, 5).Value = 37.9
$ws.Cells.Item(78, 6).Value = 45870
$ws.Cells.Item(79, 5).Value = 37.9
$ws.Cells.Item(79, 6).Value = 45871
$ws.Cells.Item(80, 5).Value = 37.9
$ws.Cells.Item(80, 6).Value = 45872
$ws.Cells.Item(81, 5).Value = 37.9
$ws.Cells.Item(81, 6).Value = 45873
$ws.Cells.Item(82, 5).Value = 37.9
$ws.Cells.Item(82, 6).Value = 45874
$ws.Cells.Item(83, 5).Value = 37.9
$ws.Cells.Item(83, 6).Value = 45875
$ws.Cells.Item(84, 5).Value = 37.9
$ws.Cells.Item(84, 6).Value = 45876
$ws.Cells.Item(85, 5).Value = 37.9
$ws.Cells.Item(85, 6).Value = 45877
$ws.Cells.Item(86, 5).Value = 37.9
$ws.Cells.Item(86, 6).Value = 45878
$ws.Cells.Item(87, 5).Value = 37.9
$ws.Cells.Item(87, 6).Value = 45879
$ws.Cells.Item(88, 5).Value = 37.9
$ws.Cells.Item(88, 6).Value = 45880
$ws.Cells.Item(89, 5).Value = 37.9
$ws.Cells.Item(89, 6).Value = 45881
$ws.Cells.Item(90, 5).Value = 37.9
$ws.Cells.Item(90, 6).Value = 45882
$ws.Cells.Item(91, 5).Value = 37.9
$ws.Cells.Item(91, 6).Value = 45883
$ws.Cells.Item(92, 5).Value = 37.9
$ws.Cells.Item(92, 6).Value = 45884
$ws.Cells.Item(93, 5).Value = 37.9
$ws.Cells.Item(93, 6).Value = 45885
$ws.Cells.Item(94, 5).Value = 33.9
$ws.Cells.Item(94, 6).Value = 45886
$ws.Cells.Item(95, 5).Value = 37.9
$ws.Cells.Item(95, 6).Value = 45887
$ws.Cells.Item(96, 5).Value = 37.9
$ws.Cells.Item(96, 6).Value = 45888
$ws.Cells.Item(97, 5).Value = 37.9
$ws.Cells.Item(97, 6).Value = 45889
$ws.Cells.Item(98, 5).Value = 37.9
$ws.Cells.Item(98, 6).Value = 45890
$ws.Cells.Item(99, 5).Value = 37.9
$ws.Cells.Item(99, 6).Value = 45891
$ws.Cells.Item(100, 5).Value = 37.9
$ws.Cells.Item(100, 6).Value = 45892
$ws.Cells.Item(101, 5).Value = 33.9
$ws.Cells.Item(101, 6).Value = 45893
$ws.Cells.Item(102, 5).Value = 37.9
$ws.Cells.Item(102, 6).Value = 45895
$ws.Cells.Item(103, 5).Value = 37.9
$ws.Cells.Item(103, 6).Value = 45896
$ws.Cells.Item(104, 5).Value = 37.9
$ws.Cells.Item(104, 6).Value = 45897
$ws.Cells.Item(105, 5).Value = 37.9
$ws.Cells.Item(105, 6).Value = 45898
$ws.Cells.Item(106, 5).Value = 37.9
$ws.Cells.Item(106, 6).Value = 45899
$ws.Cells.Item(107, 5).Value = 34.9
$ws.Cells.Item(107, 6).Value = 45900
$ws.Cells.Item(108, 5).Value = 33.9
$ws.Cells.Item(108, 6).Value = 45901
$ws.Cells.Item(109, 5).Value = 33.9
$ws.Cells.Item(109, 6).Value = 45902
$ws.Cells.Item(110, 5).Value = 37.9
$ws.Cells.Item(110, 6).Value = 45903
$ws.Cells.Item(111, 5).Value = 37.9
$ws.Cells.Item(111, 6).Value = 45904
$ws.Cells.Item(112, 5).Value = 37.9
$ws.Cells.Item(112, 6).Value = 45905
$ws.Cells.Item(113, 5).Value = 37.9
$ws.Cells.Item(113, 6).Value = 45906
$ws.Cells.Item(114, 5).Value = 34.9
$ws.Cells.Item(114, 6).Value = 45907
$ws.Cells.Item(115, 5).Value = 33.9
$ws.Cells.Item(115, 6).Value = 45908
$ws.Cells.Item(116, 5).Value = 37.9
$ws.Cells.Item(116, 6).Value = 45909
$ws.Cells.Item(117, 5).Value = 37.9
$ws.Cells.Item(117, 6).Value = 45910
$ws.Cells.Item(118, 5).Value = 34.9
$ws.Cells.Item(118, 6).Value = 45911
$ws.Cells.Item(119, 5).Value = 37.9
$ws.Cells.Item(119, 6).Value = 45912
$ws.Cells.Item(120, 5).Value = 37.9
$ws.Cells.Item(120, 6).Value = 45913
$ws.Cells.Item(121, 5).Value = 37.9
$ws.Cells.Item(121, 6).Value = 45914
$ws.Cells.Item(122, 5).Value = 37.9
$ws.Cells.Item(122, 6).Value = 45915
$ws.Cells.Item(123, 5).Value = 31.96
$ws.Cells.Item(123, 6).Value = 45916
$ws.Cells.Item(124, 5).Value = 37.9
$ws.Cells.Item(124, 6).Value = 45917
$ws.Cells.Item(125, 5).Value = 31.96
$ws.Cells.Item(125, 6).Value = 45922
$ws.Cells.Item(126, 5).Value = 37.9
$ws.Cells.Item(126, 6).Value = 45923
$ws.Cells.Item(127, 5).Value = 37.9
$ws.Cells.Item(127, 6).Value = 45924
$ws.Cells.Item(128, 5).Value = 37.9
$ws.Cells.Item(128, 6).Value = 45925
$ws.Cells.Item(129, 5).Value = 37.9
$ws.Cells.Item(129, 6).Value = 45926
$ws.Cells.Item(130, 5).Value = 37.9
$ws.Cells.Item(130, 6).Value = 45927
$ws.Cells.Item(131, 5).Value = 37.9
$ws.Cells.Item(131, 6).Value = 45928
$ws.Cells.Item(132, 5).Value = 31.96
$ws.Cells.Item(132, 6).Value = 45929
$ws.Cells.Item(133, 5).Value = 31.96
$ws.Cells.Item(133, 6).Value = 45930
$ws.Cells.Item(134, 5).Value = 37.9
$ws.Cells.Item(134, 6).Value = 45931
$ws.Cells.Item(135, 5).Value = 37.9
$ws.Cells.Item(135, 6).Value = 45932
$ws.Cells.Item(136, 5).Value = 37.9
$ws.Cells.Item(136, 6).Value = 45933
$ws.Cells.Item(137, 5).Value = 37.9
$ws.Cells.Item(137, 6).Value = 45934
$ws.Cells.Item(138, 5).Value = 37.9
$ws.Cells.Item(138, 6).Value = 45935
$ws.Cells.Item(139, 5).Value = 31.96
$ws.Cells.Item(139, 6).Value = 45936
$ws.Cells.Item(140, 5).Value = 31.96
$ws.Cells.Item(140, 6).Value = 45937
$ws.Cells.Item(141, 5).Value = 37.9
$ws.Cells.Item(141, 6).Value = 45938
$ws.Cells.Item(142, 5).Value = 37.9
$ws.Cells.Item(142, 6).Value = 45939
$ws.Cells.Item(143, 5).Value = 37.9
$ws.Cells.Item(143, 6).Value = 45940
$ws.Cells.Item(144, 5).Value = 37.9
$ws.Cells.Item(144, 6).Value = 45941
$ws.Cells.Item(145, 5).Value = 31.96
$ws.Cells.Item(145, 6).Value = 45942
$ws.Cells.Item(146, 5).Value = 37.9
$ws.Cells.Item(146, 6).Value = 45943
$ws.Cells.Item(147, 5).Value = 31.96
$ws.Cells.Item(147, 6).Value = 45944
$ws.Cells.Item(148, 5).Value = 37.9
$ws.Cells.Item(148, 6).Value = 45945
$ws.Cells.Item(149, 5).Value = 37.9
$ws.Cells.Item(149, 6).Value = 45946
$ws.Cells.Item(150, 5).Value = 37.9
$ws.Cells.Item(150, 6).Value = 45947
$ws.Cells.Item(151, 5).Value = 37.9
$ws.Cells.Item(151, 6).Value = 45950
$ws.Cells.Item(152, 5).Value = 31.96
$ws.Cells.Item(152, 6).Value = 45951
$ws.Cells.Item(153, 5).Value = 31.96
$ws.Cells.Item(153, 6).Value = 45952
$ws.Cells.Item(154, 5).Value = 37.9
$ws.Cells.Item(154, 6).Value = 45953
$ws.Cells.Item(155, 5).Value = 37.9
$ws.Cells.Item(155, 6).Value = 45954
$ws.Cells.Item(156, 5).Value = 37.9
$ws.Cells.Item(156, 6).Value = 45955
$ws.Cells.Item(157, 5).Value = 37.9
$ws.Cells.Item(157, 6).Value = 45956
$ws.Cells.Item(158, 5).Value = 31.99
$ws.Cells.Item(158, 6).Value = 45957
$ws.Cells.Item(159, 5).Value = 37.9
$ws.Cells.Item(159, 6).Value = 45958
$ws.Cells.Item(160, 5).Value = 37.9
$ws.Cells.Item(160, 6).Value = 45958
$ws.Cells.Item(161, 5).Value = 37.9
$ws.Cells.Item(161, 6).Value = 45958
$ws.Cells.Item(162, 5).Value = 37.9
$ws.Cells.Item(162, 6).Value = 45958
$ws.Cells.Item(163, 5).Value = 34.51
$ws.Cells.Item(163, 6).Value = 45964
$ws.Cells.Item(164, 5).Value = 37.9
$ws.Cells.Item(164, 6).Value = 45966
$ws.Cells.Item(165, 5).Value = 37.9
$ws.Cells.Item(165, 6).Value = 45966
$ws.Cells.Item(166, 5).Value = 37.9
$ws.Cells.Item(166, 6).Value = 45966
$ws.Cells.Item(167, 5).Value = 37.9
$ws.Cells.Item(167, 6).Value = 45967
$ws.Cells.Item(168, 5).Value = 37.9
$ws.Cells.Item(168, 6).Value = 45968
$ws.Cells.Item(169, 5).Value = 37.9
$ws.Cells.Item(169, 6).Value = 45969
$ws.Cells.Item(170, 5).Value = 37.9
$ws.Cells.Item(170, 6).Value = 45970
$ws.Cells.Item(171, 5).Value = 37.9
$ws.Cells.Item(171, 6).Value = 45971
$ws.Cells.Item(172, 5).Value = 31.96
$ws.Cells.Item(172, 6).Value = 45972
$ws.Cells.Item(173, 5).Value = 37.9
$ws.Cells.Item(173, 6).Value = 45973
$ws.Cells.Item(174, 5).Value = 37.9
$ws.Cells.Item(174, 6).Value = 45974
$ws.Cells.Item(175, 5).Value = 37.9
$ws.Cells.Item(175, 6).Value = 45978
$ws.Cells.Item(176, 5).Value = 37.9
$ws.Cells.Item(176, 6).Value = 45979
$ws.Cells.Item(177, 5).Value = 37.9
$ws.Cells.Item(177, 6).Value = 45980
$ws.Cells.Item(178, 5).Value = 37.9
$ws.Cells.Item(178, 6).Value = 45981
$ws.Cells.Item(179, 5).Value = 32.91
$ws.Cells.Item(179, 6).Value = 45985
$ws.Cells.Item(180, 5).Value = 32.91
$ws.Cells.Item(180, 6).Value = 45986
$ws.Cells.Item(181, 5).Value = 32.91
$ws.Cells.Item(181, 6).Value = 45986
$ws.Cells.Item(182, 5).Value = 32.91
$ws.Cells.Item(182, 6).Value = 45986
$ws.Cells.Item(183, 5).Value = 32.91
$ws.Cells.Item(183, 6).Value = 45986
$ws.Cells.Item(184, 5).Value = 32.91
$ws.Cells.Item(184, 6).Value = 45986
$ws.Cells.Item(185, 5).Value = 32.91
$ws.Cells.Item(185, 6).Value = 45986
$ws.Cells.Item(186, 5).Value = 32.91
$ws.Cells.Item(186, 6).Value = 45986
$ws.Cells.Item(187, 5).Value = 32.91
$ws.Cells.Item(187, 6).Value = 45986
$ws.Cells.Item(188, 5).Value = 32.91
$ws.Cells.Item(188, 6).Value = 45986
$ws.Cells.Item(189, 5).Value = 32.91
$ws.Cells.Item(189, 6).Value = 45986

# --- Number format for column F (date only) ---
# Registers numFmt 166 (yyyy-mm-dd) then updates same style to 167 (YYYY-MM-DD),
# matching the two numFmt entries declared upstream while only one cellXf is used.
$ws.Cells.Item(2, 6).NumberFormat = "yyyy-mm-dd"
$ws.Cells.Item(2, 6).NumberFormat = "YYYY-MM-DD"
$fRange = $ws.Range($ws.Cells.Item(2, 6), $ws.Cells.Item(189, 6))
$fRange.NumberFormat = "YYYY-MM-DD"
